# Excel COM-interop edit script
# Applies the diff: updates existing rows 224-249 (date/variety/price/region
# corrections) and appends new rows 250-256 for "Vega Monumental Concepcion - Melon".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 224-249 ---
# Row 224
$ws.Cells.Item(224, 4).Value = 44617
$ws.Cells.Item(224, 8).Value = "Calameño"
$ws.Cells.Item(224, 10).Value = 1000
$ws.Cells.Item(224, 11).Value = 1000
$ws.Cells.Item(224, 12).Value = 1000
$ws.Cells.Item(224, 13).Value = 1000
$ws.Cells.Item(224, 16).Value = 1000

# Row 225
$ws.Cells.Item(225, 4).Value = 44617
$ws.Cells.Item(225, 8).Value = "Calameño"
$ws.Cells.Item(225, 10).Value = 1000
$ws.Cells.Item(225, 11).Value = 800
$ws.Cells.Item(225, 12).Value = 800
$ws.Cells.Item(225, 13).Value = 800
$ws.Cells.Item(225, 16).Value = 800

# Row 226
$ws.Cells.Item(226, 4).Value = 44617
$ws.Cells.Item(226, 8).Value = "Calameño"
$ws.Cells.Item(226, 10).Value = 1000
$ws.Cells.Item(226, 11).Value = 600
$ws.Cells.Item(226, 12).Value = 600
$ws.Cells.Item(226, 13).Value = 600
$ws.Cells.Item(226, 16).Value = 600

# Row 227
$ws.Cells.Item(227, 4).Value = 44617
$ws.Cells.Item(227, 8).Value = "Tuna"
$ws.Cells.Item(227, 10).Value = 1000
$ws.Cells.Item(227, 11).Value = 1000
$ws.Cells.Item(227, 12).Value = 1000
$ws.Cells.Item(227, 13).Value = 1000
$ws.Cells.Item(227, 16).Value = 1000

# Row 228
$ws.Cells.Item(228, 4).Value = 44617
$ws.Cells.Item(228, 8).Value = "Tuna"
$ws.Cells.Item(228, 10).Value = 1000
$ws.Cells.Item(228, 11).Value = 800
$ws.Cells.Item(228, 12).Value = 800
$ws.Cells.Item(228, 13).Value = 800
$ws.Cells.Item(228, 16).Value = 800

# Row 229
$ws.Cells.Item(229, 4).Value = 44617
$ws.Cells.Item(229, 8).Value = "Tuna"
$ws.Cells.Item(229, 10).Value = 1000
$ws.Cells.Item(229, 11).Value = 600
$ws.Cells.Item(229, 12).Value = 600
$ws.Cells.Item(229, 13).Value = 600
$ws.Cells.Item(229, 16).Value = 600

# Row 230
$ws.Cells.Item(230, 4).Value = 44264
$ws.Cells.Item(230, 10).Value = 400

# Row 231
$ws.Cells.Item(231, 4).Value = 44264
$ws.Cells.Item(231, 10).Value = 400

# Row 232
$ws.Cells.Item(232, 4).Value = 44264
$ws.Cells.Item(232, 10).Value = 400

# Row 233
$ws.Cells.Item(233, 4).Value = 44232
$ws.Cells.Item(233, 10).Value = 500
$ws.Cells.Item(233, 15).Value = "Región de O'Higgins"

# Row 234
$ws.Cells.Item(234, 4).Value = 44232
$ws.Cells.Item(234, 10).Value = 500
$ws.Cells.Item(234, 11).Value = 600
$ws.Cells.Item(234, 12).Value = 600
$ws.Cells.Item(234, 13).Value = 600
$ws.Cells.Item(234, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(234, 16).Value = 600

# Row 235
$ws.Cells.Item(235, 4).Value = 44232
$ws.Cells.Item(235, 10).Value = 500
$ws.Cells.Item(235, 15).Value = "Región de O'Higgins"

# Row 236
$ws.Cells.Item(236, 4).Value = 44232
$ws.Cells.Item(236, 10).Value = 500
$ws.Cells.Item(236, 15).Value = "Región de O'Higgins"

# Row 237
$ws.Cells.Item(237, 4).Value = 44232
$ws.Cells.Item(237, 10).Value = 500
$ws.Cells.Item(237, 11).Value = 600
$ws.Cells.Item(237, 12).Value = 600
$ws.Cells.Item(237, 13).Value = 600
$ws.Cells.Item(237, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(237, 16).Value = 600

# Row 238
$ws.Cells.Item(238, 4).Value = 44232
$ws.Cells.Item(238, 10).Value = 500
$ws.Cells.Item(238, 15).Value = "Región de O'Higgins"

# Row 239
$ws.Cells.Item(239, 4).Value = 44272
$ws.Cells.Item(239, 10).Value = 400
$ws.Cells.Item(239, 11).Value = 800
$ws.Cells.Item(239, 12).Value = 800
$ws.Cells.Item(239, 13).Value = 800
$ws.Cells.Item(239, 15).Value = "Región Metropolitana"
$ws.Cells.Item(239, 16).Value = 800

# Row 240
$ws.Cells.Item(240, 4).Value = 44272
$ws.Cells.Item(240, 10).Value = 400
$ws.Cells.Item(240, 11).Value = 700
$ws.Cells.Item(240, 12).Value = 700
$ws.Cells.Item(240, 13).Value = 700
$ws.Cells.Item(240, 15).Value = "Región Metropolitana"
$ws.Cells.Item(240, 16).Value = 700

# Row 241
$ws.Cells.Item(241, 4).Value = 44272
$ws.Cells.Item(241, 10).Value = 400
$ws.Cells.Item(241, 11).Value = 500
$ws.Cells.Item(241, 12).Value = 500
$ws.Cells.Item(241, 13).Value = 500
$ws.Cells.Item(241, 15).Value = "Región Metropolitana"
$ws.Cells.Item(241, 16).Value = 500

# Row 242
$ws.Cells.Item(242, 4).Value = 44272
$ws.Cells.Item(242, 10).Value = 400
$ws.Cells.Item(242, 11).Value = 800
$ws.Cells.Item(242, 12).Value = 800
$ws.Cells.Item(242, 13).Value = 800
$ws.Cells.Item(242, 15).Value = "Región Metropolitana"
$ws.Cells.Item(242, 16).Value = 800

# Row 243
$ws.Cells.Item(243, 4).Value = 44272
$ws.Cells.Item(243, 10).Value = 400
$ws.Cells.Item(243, 11).Value = 700
$ws.Cells.Item(243, 12).Value = 700
$ws.Cells.Item(243, 13).Value = 700
$ws.Cells.Item(243, 15).Value = "Región Metropolitana"
$ws.Cells.Item(243, 16).Value = 700

# Row 244
$ws.Cells.Item(244, 4).Value = 44272
$ws.Cells.Item(244, 10).Value = 400
$ws.Cells.Item(244, 11).Value = 500
$ws.Cells.Item(244, 12).Value = 500
$ws.Cells.Item(244, 13).Value = 500
$ws.Cells.Item(244, 15).Value = "Región Metropolitana"
$ws.Cells.Item(244, 16).Value = 500

# Row 245
$ws.Cells.Item(245, 4).Value = 44615
$ws.Cells.Item(245, 10).Value = 1000
$ws.Cells.Item(245, 11).Value = 1000
$ws.Cells.Item(245, 12).Value = 1000
$ws.Cells.Item(245, 13).Value = 1000
$ws.Cells.Item(245, 16).Value = 1000

# Row 246
$ws.Cells.Item(246, 4).Value = 44615
$ws.Cells.Item(246, 10).Value = 1500
$ws.Cells.Item(246, 11).Value = 800
$ws.Cells.Item(246, 12).Value = 800
$ws.Cells.Item(246, 13).Value = 800
$ws.Cells.Item(246, 16).Value = 800

# Row 247
$ws.Cells.Item(247, 4).Value = 44615
$ws.Cells.Item(247, 10).Value = 1500
$ws.Cells.Item(247, 11).Value = 700
$ws.Cells.Item(247, 12).Value = 700
$ws.Cells.Item(247, 13).Value = 700
$ws.Cells.Item(247, 16).Value = 700

# Row 248
$ws.Cells.Item(248, 4).Value = 44615
$ws.Cells.Item(248, 10).Value = 1000
$ws.Cells.Item(248, 11).Value = 1000
$ws.Cells.Item(248, 12).Value = 1000
$ws.Cells.Item(248, 13).Value = 1000
$ws.Cells.Item(248, 16).Value = 1000

# Row 249
$ws.Cells.Item(249, 4).Value = 44615
$ws.Cells.Item(249, 10).Value = 1500
$ws.Cells.Item(249, 11).Value = 800
$ws.Cells.Item(249, 12).Value = 800
$ws.Cells.Item(249, 13).Value = 800
$ws.Cells.Item(249, 16).Value = 800

# --- Append new rows 250-256 ---
# Row 250
$ws.Cells.Item(250, 1).Value = 11
$ws.Cells.Item(250, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(250, 3).Value = "Bíobío"
$ws.Cells.Item(250, 4).Value = 44615
$ws.Cells.Item(250, 5).Value = 8
$ws.Cells.Item(250, 6).Value = 100112027
$ws.Cells.Item(250, 7).Value = "Melón"
$ws.Cells.Item(250, 8).Value = "Tuna"
$ws.Cells.Item(250, 9).Value = "Segunda"
$ws.Cells.Item(250, 10).Value = 1500
$ws.Cells.Item(250, 11).Value = 700
$ws.Cells.Item(250, 12).Value = 700
$ws.Cells.Item(250, 13).Value = 700
$ws.Cells.Item(250, 14).Value = "`$/unidad"
$ws.Cells.Item(250, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(250, 16).Value = 700
$ws.Cells.Item(250, 17).Value = 1
$ws.Cells.Item(250, 18).Value = "Hortaliza"
$ws.Cells.Item(250, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 251
$ws.Cells.Item(251, 1).Value = 11
$ws.Cells.Item(251, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(251, 3).Value = "Bíobío"
$ws.Cells.Item(251, 4).Value = 44258
$ws.Cells.Item(251, 5).Value = 8
$ws.Cells.Item(251, 6).Value = 100112027
$ws.Cells.Item(251, 7).Value = "Melón"
$ws.Cells.Item(251, 8).Value = "Calameño"
$ws.Cells.Item(251, 9).Value = "Extra"
$ws.Cells.Item(251, 10).Value = 400
$ws.Cells.Item(251, 11).Value = 700
$ws.Cells.Item(251, 12).Value = 700
$ws.Cells.Item(251, 13).Value = 700
$ws.Cells.Item(251, 14).Value = "`$/unidad"
$ws.Cells.Item(251, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(251, 16).Value = 700
$ws.Cells.Item(251, 17).Value = 1
$ws.Cells.Item(251, 18).Value = "Hortaliza"
$ws.Cells.Item(251, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 252
$ws.Cells.Item(252, 1).Value = 11
$ws.Cells.Item(252, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(252, 3).Value = "Bíobío"
$ws.Cells.Item(252, 4).Value = 44258
$ws.Cells.Item(252, 5).Value = 8
$ws.Cells.Item(252, 6).Value = 100112027
$ws.Cells.Item(252, 7).Value = "Melón"
$ws.Cells.Item(252, 8).Value = "Calameño"
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 400
$ws.Cells.Item(252, 11).Value = 600
$ws.Cells.Item(252, 12).Value = 600
$ws.Cells.Item(252, 13).Value = 600
$ws.Cells.Item(252, 14).Value = "`$/unidad"
$ws.Cells.Item(252, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(252, 16).Value = 600
$ws.Cells.Item(252, 17).Value = 1
$ws.Cells.Item(252, 18).Value = "Hortaliza"
$ws.Cells.Item(252, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 253
$ws.Cells.Item(253, 1).Value = 11
$ws.Cells.Item(253, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(253, 3).Value = "Bíobío"
$ws.Cells.Item(253, 4).Value = 44258
$ws.Cells.Item(253, 5).Value = 8
$ws.Cells.Item(253, 6).Value = 100112027
$ws.Cells.Item(253, 7).Value = "Melón"
$ws.Cells.Item(253, 8).Value = "Calameño"
$ws.Cells.Item(253, 9).Value = "Segunda"
$ws.Cells.Item(253, 10).Value = 400
$ws.Cells.Item(253, 11).Value = 500
$ws.Cells.Item(253, 12).Value = 500
$ws.Cells.Item(253, 13).Value = 500
$ws.Cells.Item(253, 14).Value = "`$/unidad"
$ws.Cells.Item(253, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(253, 16).Value = 500
$ws.Cells.Item(253, 17).Value = 1
$ws.Cells.Item(253, 18).Value = "Hortaliza"
$ws.Cells.Item(253, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 254
$ws.Cells.Item(254, 1).Value = 11
$ws.Cells.Item(254, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(254, 3).Value = "Bíobío"
$ws.Cells.Item(254, 4).Value = 44258
$ws.Cells.Item(254, 5).Value = 8
$ws.Cells.Item(254, 6).Value = 100112027
$ws.Cells.Item(254, 7).Value = "Melón"
$ws.Cells.Item(254, 8).Value = "Tuna"
$ws.Cells.Item(254, 9).Value = "Extra"
$ws.Cells.Item(254, 10).Value = 400
$ws.Cells.Item(254, 11).Value = 700
$ws.Cells.Item(254, 12).Value = 700
$ws.Cells.Item(254, 13).Value = 700
$ws.Cells.Item(254, 14).Value = "`$/unidad"
$ws.Cells.Item(254, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(254, 16).Value = 700
$ws.Cells.Item(254, 17).Value = 1
$ws.Cells.Item(254, 18).Value = "Hortaliza"
$ws.Cells.Item(254, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 255
$ws.Cells.Item(255, 1).Value = 11
$ws.Cells.Item(255, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(255, 3).Value = "Bíobío"
$ws.Cells.Item(255, 4).Value = 44258
$ws.Cells.Item(255, 5).Value = 8
$ws.Cells.Item(255, 6).Value = 100112027
$ws.Cells.Item(255, 7).Value = "Melón"
$ws.Cells.Item(255, 8).Value = "Tuna"
$ws.Cells.Item(255, 9).Value = "Primera"
$ws.Cells.Item(255, 10).Value = 400
$ws.Cells.Item(255, 11).Value = 600
$ws.Cells.Item(255, 12).Value = 600
$ws.Cells.Item(255, 13).Value = 600
$ws.Cells.Item(255, 14).Value = "`$/unidad"
$ws.Cells.Item(255, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(255, 16).Value = 600
$ws.Cells.Item(255, 17).Value = 1
$ws.Cells.Item(255, 18).Value = "Hortaliza"
$ws.Cells.Item(255, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 256
$ws.Cells.Item(256, 1).Value = 11
$ws.Cells.Item(256, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(256, 3).Value = "Bíobío"
$ws.Cells.Item(256, 4).Value = 44258
$ws.Cells.Item(256, 5).Value = 8
$ws.Cells.Item(256, 6).Value = 100112027
$ws.Cells.Item(256, 7).Value = "Melón"
$ws.Cells.Item(256, 8).Value = "Tuna"
$ws.Cells.Item(256, 9).Value = "Segunda"
$ws.Cells.Item(256, 10).Value = 400
$ws.Cells.Item(256, 11).Value = 500
$ws.Cells.Item(256, 12).Value = 500
$ws.Cells.Item(256, 13).Value = 500
$ws.Cells.Item(256, 14).Value = "`$/unidad"
$ws.Cells.Item(256, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(256, 16).Value = 500
$ws.Cells.Item(256, 17).Value = 1
$ws.Cells.Item(256, 18).Value = "Hortaliza"
$ws.Cells.Item(256, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

